$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = -356.4000064140148
$ws.Range("B2").Value = 7128

# Remove rows 3 and 4 entirely (they are no longer part of the data)
$ws.Rows("3:4").Delete()
